$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12: Aula 17 - A entidade Endereco
$ws.Cells.Item(12, 2).Value = 17
$ws.Cells.Item(12, 3).Value = "3. Classes de Domínio"
$ws.Cells.Item(12, 4).Value = "17. A entidade Endereco"

# Match style of other "observação" cells (wrap text) - copy format from E11
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("E12").Value = "3:48 - anotação @Enumerated, que define na entidade, informando para o JPA qual o tipo de dado deve ser armazenado no banco de dados com o atributo EnumType.String, salvando um enum com o tipo string"

$ws.Rows.Item(12).RowHeight = 45

# Update selection to match the new active cell state after the edit
$ws.Range("E13").Select() | Out-Null
